# Update plan-import test workbook so the SHOP 2018 QHP sheet is flagged
# as a standard plan (D2/D3 = "Y") and rename/select that tab, matching
# the "set standard_plan flag for SHOP market plans" change.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Rename "2018_QHP" -> "SHOP Q1"
$ws1.Name = "SHOP Q1"

# Flag the two plan rows as standard plans ("Yes" -> "Y")
$ws1.Range("D2").Value = "Y"
$ws1.Range("D3").Value = "Y"

# Make this sheet the active tab/selection (cell D4), matching the
# workbook/sheetView changes in the commit.
$ws1.Activate()
[void]$ws1.Range("D4").Select()
